$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C37").Value = "Y"
$ws.Range("D37").Value = "Y"
$ws.Range("C39").Value = "Y"
$ws.Range("D39").Value = "Y"
$ws.Range("C40").Value = "Y"
$ws.Range("D40").Value = "Y"
$ws.Range("C41").Value = "Y"
$ws.Range("D41").Value = "Y"
$ws.Range("C42").Value = "Y"
$ws.Range("D42").Value = "Y"
$ws.Range("C45").Value = "Y"
$ws.Range("D45").Value = "Y"
$ws.Range("C46").Value = "Y"
$ws.Range("D46").Value = "Y"
$ws.Range("C47").Value = "Y"
$ws.Range("D47").Value = "Y"
$ws.Range("C48").Value = "Y"
$ws.Range("C49").Value = "Y"
$ws.Range("D49").Value = "Y"
$ws.Range("C50").Value = "Y"
$ws.Range("D50").Value = "Y"
$ws.Range("C51").Value = "Y"
$ws.Range("D51").Value = "Y"
$ws.Range("C52").Value = "Y"
$ws.Range("D52").Value = "Y"
$ws.Range("C53").Value = "Y"
$ws.Range("D53").Value = "Y"
$ws.Range("C54").Value = "Y"
$ws.Range("D54").Value = "Y"
$ws.Range("C55").Value = "Y"
$ws.Range("D55").Value = "Y"
$ws.Range("B56").Value = "Y"
$ws.Range("C56").Value = "Y"
$ws.Range("D56").Value = "Y"
$ws.Range("B57").Value = "Y"
$ws.Range("C57").Value = "Y"
$ws.Range("D57").Value = "Y"
$ws.Range("B58").Value = "Y Segmentation fault"
$ws.Range("D58").Value = "X"
$ws.Range("B59").Value = "Y Segmentation fault"
$ws.Range("D59").Value = "X"
$ws.Range("B60").Value = "Y Segmentation fault"
$ws.Range("D60").Value = "X"
$ws.Range("B61").Value = "Y Segmentation fault"
$ws.Range("D61").Value = "X"
$ws.Range("B62").Value = "Y Segmentation fault"
$ws.Range("D62").Value = "X"
$ws.Range("B63").Value = "Y Segmentation fault"
$ws.Range("D63").Value = "X"
$ws.Range("B64").Value = "Y Segmentation fault"
$ws.Range("D64").Value = "X"
$ws.Range("B65").Value = "Y Segmentation fault"
$ws.Range("D65").Value = "X"
$ws.Range("B66").Value = "Y Segmentation fault"
$ws.Range("D66").Value = "X"
$ws.Range("B67").Value = "Y Segmentation fault"
$ws.Range("D67").Value = "X"
$ws.Range("B68").Value = "Y Segmentation fault"
$ws.Range("D68").Value = "X"
$ws.Range("B69").Value = "Y Segmentation fault"
$ws.Range("D69").Value = "X"
$ws.Range("B70").Value = "Y"
$ws.Range("C70").Value = "Y"
$ws.Range("D70").Value = "Y"
$ws.Range("B71").Value = "Y"
$ws.Range("C71").Value = "Y"
$ws.Range("D71").Value = "Y"
$ws.Range("B72").Value = "Y"
$ws.Range("C72").Value = "Y"
$ws.Range("D72").Value = "Y"
$ws.Range("B73").Value = "Y"
$ws.Range("C73").Value = "Y"
$ws.Range("D73").Value = "Y"
$ws.Range("B74").Value = "Y"
$ws.Range("C74").Value = "Y"
$ws.Range("D74").Value = "Y"
$ws.Range("B75").Value = "Y"
$ws.Range("C75").Value = "Y"
$ws.Range("D75").Value = "Y"
$ws.Range("B76").Value = "Y no Data"
$ws.Range("D76").Value = "X"
$ws.Range("B77").Value = "Y"
$ws.Range("C77").Value = "Y"
$ws.Range("D77").Value = "Y"
$ws.Range("B78").Value = "Y"
$ws.Range("C78").Value = "Y"
$ws.Range("D78").Value = "Y"
$ws.Range("B79").Value = "Y"
$ws.Range("C79").Value = "Y"
$ws.Range("D79").Value = "Y"
$ws.Range("B80").Value = "Y"
$ws.Range("C80").Value = "Y"
$ws.Range("D80").Value = "Y"
$ws.Range("B81").Value = "Y Illegal NIfTI file"
$ws.Range("D81").Value = "X"
$ws.Range("B84").Value = "Y"
$ws.Range("C84").Value = "Y"
$ws.Range("D84").Value = "Y"
$ws.Range("B85").Value = "Y"
$ws.Range("B86").Value = "Y"
$ws.Range("B87").Value = "Y"
$ws.Range("B88").Value = "Y"
$ws.Range("B89").Value = "Y"
$ws.Range("B90").Value = "Y"
$ws.Range("B91").Value = "Y"
$ws.Range("B92").Value = "Y"
$ws.Range("B93").Value = "Y"
$ws.Range("B94").Value = "Y"
$ws.Range("B95").Value = "Y"
$ws.Range("B96").Value = "Y"
$ws.Range("B97").Value = "Y"

$ws.Range("C102").Select()
$excel.ActiveWindow.ScrollRow = 70
$excel.ActiveWindow.ScrollColumn = 1
